# Add a new "Italy" test-data sheet at the end of the workbook, cloned
# from the existing "Slovakia" sheet (same layout/styles), then update
# the market name / reference number cells for Italy.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# Duplicate the Slovakia sheet and place the copy after the last sheet.
$slovakia.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Update the market name and Jira-style reference for Italy.
$italy.Range("B2").Value = "Italy Market"

$italy.Range("B4").ClearContents()
$italy.Range("B4").ClearFormats()
$italy.Range("B4").Value = "NGC-3145/T2154"

# Slovakia is no longer the active sheet/tab - deselect it and clear its
# single-cell focus in favor of a whole-sheet selection.
$slovakia.Select()
$slovakia.Cells.Select()

# Italy becomes the active sheet, with B4 focused.
$italy.Select()
$italy.Range("B4").Select()
